# ------------------------------------------------------------------
# Adds the PCA_100 / PCA_300 corpus+polarity worksheets, populates the
# PCA_30_Polarity sheet with its results table, adds the two
# "observations" notes to PCA-30-Corpus, and re-points the active tab
# to the new PCA-300-Polarity sheet.
# ------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # PCA-30-Corpus     (existing)
$ws2 = $wb.Worksheets.Item(2)          # PCA_30_Polarity   (existing, currently empty)

# --------------------------------------------------------------
# 1) PCA-30-Corpus: two new observation notes in column G
# --------------------------------------------------------------
$ws1.Range("G13").Value = "Observations : Classification is better to be performed on 1 gram than bi/tri grams"
$ws1.Range("G14").Value = "Trigram accuracy is very bad"

# Move the selection (the sheet will lose tabSelected once another
# sheet becomes active, later in the script).
$ws1.Range("E18").Select()

# --------------------------------------------------------------
# 2) PCA_30_Polarity: fill in the results table
#    Column header row + 3 ngram summary rows, re-using the exact
#    formats already present on PCA-30-Corpus (border/wrap style for
#    normal cells, bold/border style for the highlighted svm column).
# --------------------------------------------------------------

# Header row (row 1) - copy format+value from the matching PCA-30-Corpus
# header cells so the cellXfs stay de-duplicated.
$ws1.Range("A1").Copy($ws2.Range("A1")); $ws2.Range("A1").Value = "Ngram"
$ws1.Range("B1").Copy($ws2.Range("B1"))
$ws1.Range("C1").Copy($ws2.Range("C1"))
$ws1.Range("D1").Copy($ws2.Range("D1"))
$ws1.Range("E1").Copy($ws2.Range("E1"))
$ws1.Range("F1").Copy($ws2.Range("F1")); $ws2.Range("F1").Value = " svm_train"
$ws1.Range("B1").Copy($ws2.Range("G1")); $ws2.Range("G1").Value = "svm_test"
$ws1.Range("G1").Copy($ws2.Range("H1")); $ws2.Range("H1").Value = " dt_train"
$ws1.Range("H1").Copy($ws2.Range("I1"))
$ws1.Range("I1").Copy($ws2.Range("J1"))
$ws1.Range("J1").Copy($ws2.Range("K1")); $ws2.Range("K1").Value = "adaboost_test"
$ws1.Range("B1").Copy($ws2.Range("L1")); $ws2.Range("L1").Value = "rf_train"
$ws1.Range("L1").Copy($ws2.Range("M1"))
$ws2.Rows.Item(1).RowHeight = 31

# Row 2: single formatted-but-empty cell (matches the new "vertical
# center, no border / wrap" style that also appears on the other new
# sheets).
$ws2.Range("A2").Font.Size = 11
$ws2.Range("A2").VerticalAlignment = -4108
$ws2.Range("A2").WrapText = $false
$ws2.Rows.Item(2).RowHeight = 17

# Row 3 (1gram)
$ws1.Range("A2").Copy($ws2.Range("A3")); $ws2.Range("A3").Value = "1gram"
$row3 = @{ "B3"=0.60684700000000003; "C3"=0.60514100000000004; "D3"=0.52404499999999998; "E3"=0.51971999999999996;
           "F3"=0.60079499999999997; "G3"=0.59860199999999997; "H3"=0.59954300000000005; "I3"=0.59307200000000004;
           "L3"=0.60746199999999995; "M3"=0.60267099999999996 }
foreach ($addr in $row3.Keys) {
    $ws1.Range("B2").Copy($ws2.Range($addr))
    $ws2.Range($addr).Value = $row3[$addr]
}
$ws1.Range("J2").Copy($ws2.Range("J3")); $ws2.Range("J3").Value = 0.69508700000000001
$ws1.Range("K2").Copy($ws2.Range("K3")); $ws2.Range("K3").Value = 0.65692799999999996
$ws2.Rows.Item(3).RowHeight = 17

# Row 4 (2gram)
$ws1.Range("A2").Copy($ws2.Range("A4")); $ws2.Range("A4").Value = "2gram"
$row4 = @{ "B4"=0.53445600000000004; "C4"=0.53599699999999995; "D4"=0.50204000000000004; "E4"=0.50560000000000005;
           "F4"=0.528277; "G4"=0.53102400000000005; "H4"=0.56964300000000001; "I4"=0.56406500000000004;
           "J4"=0.62340300000000004; "K4"=0.59536699999999998; "L4"=0.56852999999999998; "M4"=0.56027400000000005 }
foreach ($addr in $row4.Keys) {
    $ws1.Range("B2").Copy($ws2.Range($addr))
    $ws2.Range($addr).Value = $row4[$addr]
}
$ws2.Rows.Item(4).RowHeight = 17

# Row 5 (3gram)
$ws1.Range("A2").Copy($ws2.Range("A5")); $ws2.Range("A5").Value = "3gram"
$row5 = @{ "B5"=0.49860900000000002; "C5"=0.49867800000000001; "D5"=0.49233700000000002; "E5"=0.49179200000000001;
           "F5"=0.498446; "G5"=0.49867800000000001; "H5"=0.52354699999999998; "I5"=0.52222500000000005;
           "J5"=0.53935999999999995; "K5"=0.52876299999999998; "L5"=0.52179600000000004; "M5"=0.52187700000000004 }
foreach ($addr in $row5.Keys) {
    $ws1.Range("B2").Copy($ws2.Range($addr))
    $ws2.Range($addr).Value = $row5[$addr]
}
$ws2.Rows.Item(5).RowHeight = 17

$ws2.Range("K3").Select()

# --------------------------------------------------------------
# 3) Add the four new worksheets, in tab order:
#    PCA_100_corpus, PCA_100_Polarity, PCA-300-Polarity, PCA-300-Corpus
# --------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PCA_100_corpus"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "PCA_100_Polarity"

# --------------------------------------------------------------
# 4) PCA_100_Polarity: same shape as PCA_30_Polarity but a 3-row
#    table (no blank row2) and plain (unstyled) row labels.
# --------------------------------------------------------------
$ws4.Range("A1").Value = "ngrams"
$ws1.Range("B1").Copy($ws4.Range("B1"))
$ws1.Range("C1").Copy($ws4.Range("C1")); $ws4.Range("C1").Value = "logreg_test"
$ws1.Range("D1").Copy($ws4.Range("D1"))
$ws1.Range("E1").Copy($ws4.Range("E1"))
$ws1.Range("F1").Copy($ws4.Range("F1")); $ws4.Range("F1").Value = " svm_train"
$ws1.Range("G1").Copy($ws4.Range("G1")); $ws4.Range("G1").Value = " svm_test"
$ws1.Range("G1").Copy($ws4.Range("H1")); $ws4.Range("H1").Value = " dt_train"
$ws1.Range("H1").Copy($ws4.Range("I1"))
$ws1.Range("I1").Copy($ws4.Range("J1"))
$ws1.Range("J1").Copy($ws4.Range("K1")); $ws4.Range("K1").Value = " adaboost_test"
$ws1.Range("L1").Copy($ws4.Range("L1"))
$ws1.Range("L1").Copy($ws4.Range("M1")); $ws4.Range("M1").Value = " rf_test"
$ws4.Rows.Item(1).RowHeight = 46

$ws4.Range("A2").Value = "1gram"
$row2b = @{ "B2"=0.67854199999999998; "C2"=0.68189999999999995; "D2"=0.56522600000000001; "E2"=0.56048299999999995;
            "F2"=0.67477500000000001; "G2"=0.67817899999999998; "H2"=0.62526999999999999; "I2"=0.61794000000000004;
            "L2"=0.65472900000000001; "M2"=0.65358899999999998 }
foreach ($addr in $row2b.Keys) {
    $ws1.Range("B2").Copy($ws4.Range($addr))
    $ws4.Range($addr).Value = $row2b[$addr]
}
$ws1.Range("J2").Copy($ws4.Range("J2")); $ws4.Range("J2").Value = 0.743838
$ws1.Range("K2").Copy($ws4.Range("K2")); $ws4.Range("K2").Value = 0.71035099999999995
$ws4.Rows.Item(2).RowHeight = 17

$ws4.Range("A3").Value = "2gram"
$row3b = @{ "B3"=0.57621699999999998; "C3"=0.575021; "D3"=0.53829400000000005; "E3"=0.53575399999999995;
            "F3"=0.57169499999999995; "G3"=0.57199500000000003; "H3"=0.57060500000000003; "I3"=0.56566499999999997;
            "J3"=0.64034100000000005; "K3"=0.60434100000000002; "L3"=0.58527099999999999; "M3"=0.57481199999999999 }
foreach ($addr in $row3b.Keys) {
    $ws1.Range("B2").Copy($ws4.Range($addr))
    $ws4.Range($addr).Value = $row3b[$addr]
}
$ws4.Rows.Item(3).RowHeight = 17

$ws4.Range("A4").Value = "3gram"
$row4b = @{ "B4"=0.51828300000000005; "C4"=0.51346000000000003; "D4"=0.50244599999999995; "E4"=0.49627900000000003;
            "F4"=0.51783100000000004; "G4"=0.512521; "H4"=0.51711200000000002; "I4"=0.51196399999999997;
            "J4"=0.54228200000000004; "K4"=0.52806799999999998; "L4"=0.52515800000000001; "M4"=0.51999899999999999 }
foreach ($addr in $row4b.Keys) {
    $ws1.Range("B2").Copy($ws4.Range($addr))
    $ws4.Range($addr).Value = $row4b[$addr]
}
$ws4.Rows.Item(4).RowHeight = 17

$ws4.Range("M11").Select()

# --------------------------------------------------------------
# 5) PCA-300-Polarity: plain (unstyled) ngram-vs-corpus listing
# --------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "PCA-300-Polarity"
$ws5.Range("A1").Value = "ngram"
$ws5.Range("A2").Value = "1gram"
$ws5.Range("A3").Value = "2gram"
$ws5.Range("A4").Value = "3gram"

# --------------------------------------------------------------
# 6) PCA-300-Corpus: new, still empty
# --------------------------------------------------------------
$ws6 = $wb.Worksheets.Add($null, $ws5)
$ws6.Name = "PCA-300-Corpus"

# --------------------------------------------------------------
# 7) Final UI state: PCA-300-Polarity (tab index 5 / activeTab=4) is
#    the active sheet, with A4 selected.
# --------------------------------------------------------------
$ws5.Activate()
$ws5.Range("A4").Select()
